$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet "E2E_002" right after "E2E_001"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "E2E_002"

# ---- Header row ----
$ws2.Range("A1").Value = "SlNo"
$ws2.Range("B1").Value = "Description"
$ws2.Range("C1").Value = "Action"
$ws2.Range("D1").Value = "Element Reference"
$ws2.Range("E1").Value = "Reference Value "
$ws2.Range("F1").Value = "Data"
$ws2.Range("G1").Value = "Result"

# copy header formatting from sheet1 header row
$ws1.Range("A1:G1").Copy() | Out-Null
$ws2.Range("A1:G1").PasteSpecial(-4122) | Out-Null

# ---- Data rows 2-6 (first test run) ----
$data = @(
    @{ A=1; B="Launch MMP Application"; C="launchApplication"; D="null"; E="http://96.84.175.78/MMP-Release2-Integrated-Build.6.8.000/portal/login.php"; F="null" },
    @{ A=2; B="Enter Username"; C="findElement"; D="id"; E="username"; F="ria12345" },
    @{ A=3; B="Enter Password"; C="findElement"; D="id"; E="password"; F="Ria12345" },
    @{ A=4; B="Submit Button"; C="findElement"; D="name"; E="submit"; F="null" },
    @{ A=5; B="Verifying invalid login Message"; C="verifyTextInAlert"; D="null"; E="null"; F="Wrong username and password."; G="Pass" }
)

# ---- Data rows 7-10 (second test run, mirrors rows 2-5) ----
$data2 = @(
    @{ A=1; B="Launch MMP Application"; C="launchApplication"; D="null"; E="http://96.84.175.78/MMP-Release2-Integrated-Build.6.8.000/portal/login.php"; F="null" },
    @{ A=2; B="Enter Username"; C="findElement"; D="id"; E="username"; F="ria1" },
    @{ A=3; B="Enter Password"; C="findElement"; D="id"; E="password"; F="Ria12345" },
    @{ A=4; B="Submit Button"; C="findElement"; D="name"; E="submit"; F="null" }
)

# copy the row-format template from sheet1 (rows 2-6 cover every style variant we need)
$ws1.Range("A2:G2").Copy() | Out-Null
$ws2.Range("A2:G2").PasteSpecial(-4122) | Out-Null
$ws1.Range("A3:G3").Copy() | Out-Null
$ws2.Range("A3:G3").PasteSpecial(-4122) | Out-Null
$ws1.Range("A4:G4").Copy() | Out-Null
$ws2.Range("A4:G4").PasteSpecial(-4122) | Out-Null
$ws1.Range("A3:G3").Copy() | Out-Null
$ws2.Range("A5:G5").PasteSpecial(-4122) | Out-Null
$ws1.Range("A7:G7").Copy() | Out-Null
$ws2.Range("A6:G6").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2:G2").Copy() | Out-Null
$ws2.Range("A7:G7").PasteSpecial(-4122) | Out-Null
$ws1.Range("A3:G3").Copy() | Out-Null
$ws2.Range("A8:G8").PasteSpecial(-4122) | Out-Null
$ws1.Range("A4:G4").Copy() | Out-Null
$ws2.Range("A9:G9").PasteSpecial(-4122) | Out-Null
$ws1.Range("A3:G3").Copy() | Out-Null
$ws2.Range("A10:G10").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws2.Cells.Item($r, 1).Value = $row.A
    $ws2.Cells.Item($r, 2).Value = $row.B
    $ws2.Cells.Item($r, 3).Value = $row.C
    $ws2.Cells.Item($r, 4).Value = $row.D
    $ws2.Cells.Item($r, 5).Value = $row.E
    $ws2.Cells.Item($r, 6).Value = $row.F
    if ($row.ContainsKey("G")) {
        $ws2.Cells.Item($r, 7).Value = $row.G
    }
}

for ($i = 0; $i -lt $data2.Count; $i++) {
    $r = 7 + $i
    $row = $data2[$i]
    $ws2.Cells.Item($r, 1).Value = $row.A
    $ws2.Cells.Item($r, 2).Value = $row.B
    $ws2.Cells.Item($r, 3).Value = $row.C
    $ws2.Cells.Item($r, 4).Value = $row.D
    $ws2.Cells.Item($r, 5).Value = $row.E
    $ws2.Cells.Item($r, 6).Value = $row.F
}

# Hyperlinks for E2 and E7 (launch URL cells)
$ws2.Hyperlinks.Add($ws2.Range("E2"), "http://96.84.175.78/MMP-Release2-Integrated-Build.6.8.000/portal/login.php") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E7"), "http://96.84.175.78/MMP-Release2-Integrated-Build.6.8.000/portal/login.php") | Out-Null

# ---- column widths ----
$ws2.Columns.Item(2).ColumnWidth = 16.54296875
$ws2.Columns.Item(3).ColumnWidth = 15
$ws2.Columns.Item(4).ColumnWidth = 19.36328125
$ws2.Columns.Item(5).ColumnWidth = 36.7265625
$ws2.Columns.Item(6).ColumnWidth = 11.81640625

# ---- sheet view ----
$ws2.Range("E9").Select() | Out-Null

# Set the active sheet to E2E_002 (it becomes the new active/selected tab)
$ws2.Activate() | Out-Null
